# Apply weekly fruit/vegetable price updates (permuted row values) per commit "Fruta / hortaliza, semanal"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44215
$ws.Range("J2").Value = 60
$ws.Range("K2").Value = 35000
$ws.Range("L2").Value = 35000
$ws.Range("M2").Value = 35000
$ws.Range("P2").Value = 1400

# Row 4
$ws.Range("D4").Value = 44239
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 35000
$ws.Range("L4").Value = 35000
$ws.Range("M4").Value = 35000
$ws.Range("P4").Value = 1400

# Row 5
$ws.Range("D5").Value = 44218
$ws.Range("K5").Value = 42000
$ws.Range("L5").Value = 42000
$ws.Range("M5").Value = 42000
$ws.Range("O5").Value = "Región del Maule"
$ws.Range("P5").Value = 1680

# Row 6
$ws.Range("D6").Value = 44243
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 33000
$ws.Range("L6").Value = 33000
$ws.Range("M6").Value = 33000
$ws.Range("P6").Value = 1320

# Row 7
$ws.Range("D7").Value = 44250
$ws.Range("J7").Value = 70

# Row 8
$ws.Range("D8").Value = 44253
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 30000
$ws.Range("L8").Value = 30000
$ws.Range("M8").Value = 30000
$ws.Range("P8").Value = 1200

# Row 9
$ws.Range("D9").Value = 44211
$ws.Range("K9").Value = 42000
$ws.Range("L9").Value = 42000
$ws.Range("M9").Value = 42000
$ws.Range("P9").Value = 1680

# Row 11
$ws.Range("D11").Value = 44203
$ws.Range("J11").Value = 20
$ws.Range("K11").Value = 30000
$ws.Range("L11").Value = 30000
$ws.Range("M11").Value = 30000
$ws.Range("P11").Value = 1200

# Row 12
$ws.Range("D12").Value = 44586
$ws.Range("J12").Value = 80
$ws.Range("K12").Value = 31000
$ws.Range("L12").Value = 31000
$ws.Range("M12").Value = 31000
$ws.Range("O12").Value = "Región del Maule"
$ws.Range("P12").Value = 1240

# Row 13
$ws.Range("D13").Value = 44246
$ws.Range("I13").Value = "Primera"
$ws.Range("K13").Value = 31000
$ws.Range("L13").Value = 31000
$ws.Range("M13").Value = 31000
$ws.Range("O13").Value = "Región del Maule"
$ws.Range("P13").Value = 1240

# Row 14
$ws.Range("D14").Value = 44589
$ws.Range("J14").Value = 90
$ws.Range("K14").Value = 31000
$ws.Range("L14").Value = 31000
$ws.Range("M14").Value = 31000
$ws.Range("O14").Value = "Región Metropolitana"
$ws.Range("P14").Value = 1240

# Row 15
$ws.Range("D15").Value = 44582
$ws.Range("J15").Value = 40
$ws.Range("K15").Value = 35000
$ws.Range("L15").Value = 35000
$ws.Range("M15").Value = 35000
$ws.Range("O15").Value = "Región Metropolitana"
$ws.Range("P15").Value = 1400

# Row 16
$ws.Range("D16").Value = 44582
$ws.Range("I16").Value = "Segunda"
$ws.Range("J16").Value = 40
$ws.Range("K16").Value = 27000
$ws.Range("L16").Value = 27000
$ws.Range("M16").Value = 27000
$ws.Range("P16").Value = 1080

# Row 17
$ws.Range("D17").Value = 44202
$ws.Range("J17").Value = 30
$ws.Range("K17").Value = 30000
$ws.Range("L17").Value = 30000
$ws.Range("M17").Value = 30000
$ws.Range("O17").Value = "Región del Maule"
$ws.Range("P17").Value = 1200

# Row 18
$ws.Range("D18").Value = 44204
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 50
$ws.Range("K18").Value = 37000
$ws.Range("L18").Value = 37000
$ws.Range("M18").Value = 37000
$ws.Range("O18").Value = "Región del Maule"
$ws.Range("P18").Value = 1480

# Row 19
$ws.Range("D19").Value = 44201
$ws.Range("J19").Value = 60
$ws.Range("K19").Value = 30000
$ws.Range("L19").Value = 30000
$ws.Range("M19").Value = 30000
$ws.Range("O19").Value = "Región del Maule"
$ws.Range("P19").Value = 1200

# Row 20
$ws.Range("D20").Value = 44579
$ws.Range("K20").Value = 35000
$ws.Range("L20").Value = 35000
$ws.Range("M20").Value = 35000
$ws.Range("O20").Value = "Región Metropolitana"
$ws.Range("P20").Value = 1400

# Row 21
$ws.Range("D21").Value = 44579
$ws.Range("I21").Value = "Segunda"
$ws.Range("J21").Value = 50
$ws.Range("K21").Value = 30000
$ws.Range("L21").Value = 30000
$ws.Range("M21").Value = 30000
$ws.Range("O21").Value = "Región Metropolitana"
$ws.Range("P21").Value = 1200

# Row 22
$ws.Range("D22").Value = 44225
$ws.Range("J22").Value = 60
$ws.Range("K22").Value = 32000
$ws.Range("L22").Value = 32000
$ws.Range("M22").Value = 32000
$ws.Range("P22").Value = 1280

# Row 23
$ws.Range("D23").Value = 44271
$ws.Range("J23").Value = 40

# Row 24
$ws.Range("D24").Value = 44568
$ws.Range("J24").Value = 60
$ws.Range("K24").Value = 35000
$ws.Range("L24").Value = 35000
$ws.Range("M24").Value = 35000
$ws.Range("O24").Value = "Región Metropolitana"
$ws.Range("P24").Value = 1400

# Row 25
$ws.Range("D25").Value = 44264
$ws.Range("K25").Value = 29000
$ws.Range("L25").Value = 29000
$ws.Range("M25").Value = 29000
$ws.Range("O25").Value = "Región Metropolitana"
$ws.Range("P25").Value = 1160

# Row 26
$ws.Range("D26").Value = 44232
$ws.Range("J26").Value = 40
$ws.Range("K26").Value = 40000
$ws.Range("L26").Value = 40000
$ws.Range("M26").Value = 40000
$ws.Range("P26").Value = 1600

# Row 27
$ws.Range("D27").Value = 44572
$ws.Range("J27").Value = 80
$ws.Range("K27").Value = 35000
$ws.Range("L27").Value = 35000
$ws.Range("M27").Value = 35000
$ws.Range("O27").Value = "Región Metropolitana"
$ws.Range("P27").Value = 1400

# Row 28
$ws.Range("D28").Value = 44236
$ws.Range("J28").Value = 40
$ws.Range("K28").Value = 38000
$ws.Range("L28").Value = 38000
$ws.Range("M28").Value = 38000
$ws.Range("O28").Value = "Región del Maule"
$ws.Range("P28").Value = 1520
